$d = $word.ActiveDocument

# --- 1. Remove the "_h2o_keep_element" wrapper paragraphs -----------------
# These are the NodeStart / HeadSeparator / HeadEnd / NodeEnd paragraphs
# that used to bracket the resource content. We locate them by style name
# (rather than fixed index) so the script is robust to ordering.

$stylesToStrip = @("NodeStart", "HeadSeparator", "HeadEnd", "NodeEnd")

foreach ($styleName in $stylesToStrip) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq $d.Styles.Item($styleName).NameLocal) {
            $p.Range.Delete()
        }
    }
}

# --- 2. Insert a new "invisibleseparator" paragraph before the horizontal
#        rule paragraph (the empty paragraph that holds the <w:pict> rule).
#        We add it *after* the last body paragraph (rather than *before* the
#        rule paragraph) which keeps Word's split-point bookkeeping minimal.

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match [char]13 -and $p.Range.Text.Trim([char]13).Length -eq 0) {
        $prev = $d.Paragraphs.Item($i - 1)
        $prev.Range.InsertParagraphAfter()
        $sep = $d.Paragraphs.Item($i)
        $sep.Style = "invisibleseparator"
        $sep.Range.Text = " "
        break
    }
}

# --- 3. Remove the now-unused custom styles --------------------------------

foreach ($styleName in $stylesToStrip) {
    $style = $d.Styles.Item($styleName)
    $style.Delete()
}
